# "Add files via upload" — the author re-uploaded destinatarios.xlsx with the
# recipient list simplified: the two "nombre destinatario" entries were
# replaced with placeholder names and the per-row WhatsApp comment text
# (column C) was cleared out for both data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Hans Fendt -> Hans 1, drop the comment text in C2
$ws.Range("B2").Value = "Hans 1"
$ws.Range("C2").ClearContents()

# Row 3: Eduardo Riera -> Hans 2, drop the comment text in C3
$ws.Range("B3").Value = "Hans 2"
$ws.Range("C3").ClearContents()

# Match the saved selection in the refreshed file (active cell B4)
$ws.Range("B4").Select() | Out-Null
